{"js": "// Office.js (Word JavaScript API) script.\n// Appends \"O que aprendemos: \" to the trailing (empty/space) paragraph of\n// the document and adds three new bulleted (\"O que aprendemos\") list items\n// describing what was learned in the lesson, matching the commit\n// \"Finalizei a aula 1 / o que aprendemos\".\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// The last paragraph in the document is the trailing list item that only\n// contains a single space - this is where \"O que aprendemos: \" gets\n// appended.\nconst items = paragraphs.items;\nconst lastParagraph = items[items.length - 1];\n\n// Append the new text run to the end of the existing (last) paragraph,\n// instead of creating a new paragraph.\nconst endRange = lastParagraph.getRange(Word.RangeLocation.end);\nendRange.insertText(\"O que aprendemos: \", Word.InsertLocation.replace);\nawait context.sync();\n\n// Helper that inserts a new paragraph after `afterParagraph`, with the\n// given text, cloning the preceding paragraph's formatting (style,\n// numbering, spacing, run fonts, ...) and then forcing the list level\n// (ilvl) to 2 (0-based) to match the sub-bullets used elsewhere in this\n// document for the \"O que aprendemos\" summaries.\nfunction addSubBullet(afterParagraph, text) {\n  const newParagraph = afterParagraph.insertParagraph(text, Word.InsertLocation.after);\n  const listItem = newParagraph.listItemOrNullObject;\n  listItem.level = 2;\n  return newParagraph;\n}\n\nlet anchor = lastParagraph;\nanchor = addSubBullet(\n  anchor,\n  \"URLs e seus formatos: como as URLs funcionam e o que cada parte de uma URL significa - base e par\u00e2metros;\"\n);\nanchor = addSubBullet(\n  anchor,\n  \"O operador de fatiamento [a:b], utilizado para obter uma substring desde o \u00edndice a at\u00e9 o \u00edndice b - 1 da string original. Lembrando que b - 1 pois o segundo argumento do fatiamento \u00e9 exclusivo;\"\n);\nanchor = addSubBullet(\n  anchor,\n  \"A string original n\u00e3o \u00e9 alterada ao ser fatiada devido \u00e0 sua imutabilidade.\"\n);\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# Appends \"O que aprendemos: \" to the trailing (empty/space) paragraph of\n# the document and adds three new bulleted (\"O que aprendemos\") list items\n# describing what was learned in the lesson, matching the commit\n# \"Finalizei a aula 1 / o que aprendemos\".\n\n$d = $word.ActiveDocument\n\n# The last paragraph in the document is the trailing list item that only\n# contains a single space - this is where \"O que aprendemos: \" gets\n# appended (rather than creating a brand new paragraph).\n$count = $d.Paragraphs.Count\n$lastParagraph = $d.Paragraphs.Item($count)\n$endRange = $lastParagraph.Range\n$endRange.Collapse(0)   # wdCollapseEnd\n$endRange.InsertAfter(\"O que aprendemos: \")\n\n# Texts for the three new sub-bullets summarizing what was learned.\n$bullets = @(\n    \"URLs e seus formatos: como as URLs funcionam e o que cada parte de uma URL significa - base e par\u00e2metros;\",\n    \"O operador de fatiamento [a:b], utilizado para obter uma substring desde o \u00edndice a at\u00e9 o \u00edndice b - 1 da string original. Lembrando que b - 1 pois o segundo argumento do fatiamento \u00e9 exclusivo;\",\n    \"A string original n\u00e3o \u00e9 alterada ao ser fatiada devido \u00e0 sua imutabilidade.\"\n)\n\n# Insert each bullet as a new paragraph right after the previous one,\n# cloning the paragraph formatting (style/numbering/spacing/fonts) from the\n# paragraph it follows, then force the list level (ListLevelNumber, which\n# is 1-based) to 3 so the OOXML ilvl ends up at 2 - matching the other\n# third-level sub-bullets (\"O fatiamento de strings...\", etc.) already used\n# throughout this document.\n$anchorRange = $lastParagraph.Range\n$anchorRange.Collapse(0)   # wdCollapseEnd\n\nforeach ($bulletText in $bullets) {\n    $anchorRange.InsertParagraphAfter()\n    $newCount = $d.Paragraphs.Count\n    $newParagraph = $d.Paragraphs.Item($newCount)\n    $newParagraph.Range.Text = $bulletText\n    $newParagraph.Range.ListFormat.ListLevelNumber = 3\n    $anchorRange = $newParagraph.Range\n    $anchorRange.Collapse(0)   # wdCollapseEnd\n}\n"}
